# New datasets + baseline regression
# Rename the "congenital" category to "misc_long_term" across every
# variables_* worksheet in the workbook (it shows up in column A, on
# either row 3 or row 4 depending on the sheet).

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    $lastRow = $used.Rows.Count

    for ($r = 1; $r -le $lastRow; $r++) {
        $cell = $ws.Cells.Item($r, 1)
        if ($cell.Value2 -eq "congenital") {
            $cell.Value = "misc_long_term"
        }
    }
}
